$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. The worksheet that used to be "总计" (sheetId 6) becomes "2022-Q1"
#    and is repopulated with the new fund-holdings table.
# ------------------------------------------------------------------
$q1 = $wb.Worksheets.Item(6)
$q1.Name = "2022-Q1"

# Header row (B1 already carries the bold/centered/bordered style used
# throughout the workbook - reuse it and extend it across E1:H1).
$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("B1").Copy()
$q1.Range("E1:H1").PasteSpecial(-4122)
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

# Row-index column (A) keeps the same bold/centered/bordered style as
# the existing A2:A6 cells; extend that style down to A7:A12.
$q1.Range("A2").Copy()
$q1.Range("A7:A12").PasteSpecial(-4122)
for ($i = 0; $i -lt 11; $i++) {
    $q1.Cells.Item($i + 2, 1).Value = $i
}

# Make columns B and D:G text cells (matches the source data which
# stores these numbers as plain text) so values like "006424" or
# "90.90" are not silently turned into numbers / stripped of leading
# zeros.
$q1.Range("B2:B12").NumberFormat = "@"
$q1.Range("D2:G12").NumberFormat = "@"

$funds = @(
    @("202027", "南方高端装备灵活配置混合A", "15.67", "90.90", "3.67", "0.5751", 6),
    @("010452", "广发瑞福精选混合A", "16.29", "78.69", "2.97", "0.4838", 7),
    @("006424", "嘉合锦程价值精选混合A", "3.45", "82.01", "6.64", "0.2291", 3),
    @("012987", "嘉合锦明混合A", "6.24", "63.80", "3.54", "0.2209", 5),
    @("012988", "嘉合锦明混合C", "3.82", "63.80", "3.54", "0.1352", 5),
    @("006425", "嘉合锦程价值精选混合C", "1.35", "82.01", "6.64", "0.0896", 3),
    @("005207", "南方高端装备灵活配置混合C", "2.06", "90.90", "3.67", "0.0756", 6),
    @("010453", "广发瑞福精选混合C", "2.42", "78.69", "2.97", "0.0719", 7),
    @("007141", "嘉合稳健增长灵活配置混合A", "0.83", "76.90", "3.98", "0.0330", 10),
    @("007142", "嘉合稳健增长灵活配置混合C", "0.50", "76.90", "3.98", "0.0199", 10),
    @("007281", "嘉合消费升级混合", "0.24", "82.28", "4.80", "0.0115", 9)
)

$r = 2
foreach ($fund in $funds) {
    $q1.Cells.Item($r, 2).Value = $fund[0]
    $q1.Cells.Item($r, 3).Value = $fund[1]
    $q1.Cells.Item($r, 4).Value = $fund[2]
    $q1.Cells.Item($r, 5).Value = $fund[3]
    $q1.Cells.Item($r, 6).Value = $fund[4]
    $q1.Cells.Item($r, 7).Value = $fund[5]
    $q1.Cells.Item($r, 8).Value = $fund[6]
    $r = $r + 1
}

# ------------------------------------------------------------------
# 2. A brand-new worksheet named "总计" is inserted right after
#    "2022-Q1" and holds the refreshed summary table (the previous
#    "总计" rows, shifted down one, plus the new 2022-Q1 row on top).
# ------------------------------------------------------------------
$zj = $wb.Worksheets.Add($null, $q1)
$zj.Name = "总计"

$zj.Range("B1").Value = "日期"
$zj.Range("C1").Value = "持有数量(只)"
$zj.Range("D1").Value = "持有市值(亿元)"
$q1.Range("B1").Copy()
$zj.Range("B1:D1").PasteSpecial(-4122)

$q1.Range("A2").Copy()
$zj.Range("A2:A7").PasteSpecial(-4122)

$summary = @(
    @("2022-Q1", 11, 1.95),
    @("2021-Q4", 4, 0.2),
    @("2021-Q3", 75, 23.21),
    @("2021-Q2", 22, 4.8),
    @("2021-Q1", 11, 1.04),
    @("2020-Q4", 19, 1.05)
)

$r = 2
$idx = 0
foreach ($row in $summary) {
    $zj.Cells.Item($r, 1).Value = $idx
    $zj.Cells.Item($r, 2).Value = $row[0]
    $zj.Cells.Item($r, 3).Value = $row[1]
    $zj.Cells.Item($r, 4).Value = $row[2]
    $r = $r + 1
    $idx = $idx + 1
}
